$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: refine the timestamp's stored precision (same instant, tighter float)
$ws.Range("A16").Value2 = 45877.62522983796

# New row 17: append the next weather-station reading
$ws.Range("A17").NumberFormat = $ws.Range("A16").NumberFormat
$ws.Range("A17").Value2 = 45877.70853273479
$ws.Range("B17").Value2 = 2025
$ws.Range("C17").Value2 = 32
$ws.Range("D17").Value2 = 19.17
$ws.Range("E17").Value2 = 76.87
$ws.Range("F17").Value2 = 138.87
$ws.Range("G17").Value2 = 10.93
$ws.Range("H17").Value = "ESE"
$ws.Range("I17").Value2 = 0
$ws.Range("J17").Value = "17:00:17"
